$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top, shifting existing rows down
$ws.Rows.Item(1).Insert()

# Set the new cell's value
$ws.Range("A1").Value = "text"

# Move selection to A2, matching typical Excel behavior after entering a value
$ws.Range("A2").Select()
